$wb = $excel.ActiveWorkbook

# Update the "展览" sheet (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 99
$ws1.Range("F6").Value = 7

# Update the "全部类型" sheet (All types) which mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 99
$ws4.Range("F6").Value = 7
